# Update the two worksheets so that each gains a "positives" column (B)
# and a "pos_ranked" column (D), with the original "ranks" column shifting
# to column C. Some row labels (and their ranks) also get reordered.

$wb = $excel.ActiveWorkbook

function Set-SheetData {
    param(
        [string]$SheetName,
        [object[][]]$Rows   # each row: name, positives, ranks, pos_ranked
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Header row
    $ws.Cells.Item(1, 1).Value = "name"
    $ws.Cells.Item(1, 2).Value = "positives"
    $ws.Cells.Item(1, 3).Value = "ranks"
    $ws.Cells.Item(1, 4).Value = "pos_ranked"

    # Make sure the new header cells (C1, D1) carry the same formatting as
    # the existing header cells (bold + centered), matching A1/B1.
    $ws.Range("C1:D1").Font.Bold = $true
    $ws.Range("C1:D1").HorizontalAlignment = -4108

    $r = 2
    foreach ($row in $Rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $r++
    }
}

$sheet1Rows = @(
    ,@("state", 100, 100, 100)
    ,@("region", 94, 87, 87)
    ,@("rplthemes", 66, 54, 54)
    ,@("percentblackorafricanamerican", 57.99999999999999, 44, 44)
    ,@("derivedtotalenrolled", 50, 42, 42)
    ,@("percentstudentsfreereducedlunch", 59, 41, 41)
    ,@("percenttwoormoreraces", 40, 27, 27)
    ,@("percentamericanindianoralaskanative", 45, 24, 24)
    ,@("cntycaseschange", 47, 17, 17)
    ,@("locale", 32, 18, 17)
    ,@("percenthispaniclatino", 35, 12, 12)
    ,@("percentnativehawaiianorotherpacificislander", 24, 12, 12)
    ,@("percentwhite", 33, 12, 12)
    ,@("percentasian", 29, 6, 6)
    ,@("percentnotspecified", 8, 4, 0)
)

$sheet2Rows = @(
    ,@("vaccination", 100, 100, 100)
    ,@("testingandscreening", 100, 99, 99)
    ,@("physicaldistancing", 49, 51, 47)
    ,@("etiquette", 38, 32, 37)
    ,@("traceandquarantine", 30, 25, 29)
    ,@("masks", 27, 24, 27)
    ,@("ventilation", 21, 16, 20)
    ,@("cleaning", 8, 6, 7.000000000000001)
    ,@("cohortingorstaggering", 0, 47, 0)
    ,@("stayhome", 0, 0, 0)
)

Set-SheetData "covariate_importance" $sheet1Rows
Set-SheetData "strategy_importance" $sheet2Rows
